# Update "想去人数" (want-to-go count) figures in the 展览 (sheet 1) and
# 全部类型 (sheet 4) worksheets to match the freshly generated gh-pages data.

$wb = $excel.ActiveWorkbook

# Rows whose F-column value changes identically on both sheets.
$commonUpdates = @{
    5  = 51
    6  = 185
    8  = 44
    10 = 16
    11 = 42
    14 = 1539
    16 = 488
    17 = 450
    18 = 141
    22 = 1407
    23 = 3341
    26 = 56
    27 = 1082
    29 = 1774
    30 = 565
    34 = 399
    36 = 637
    38 = 31
}

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)
    foreach ($row in $commonUpdates.Keys) {
        $ws.Cells.Item($row, 6).Value = $commonUpdates[$row]
    }
}

# Row 32 ends at the same value (52) on both sheets but started from a
# different number on each, so it's applied separately.
$wb.Worksheets.Item(1).Cells.Item(32, 6).Value = 52
$wb.Worksheets.Item(4).Cells.Item(32, 6).Value = 52
